$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so numeric-looking strings
# (e.g. "602.15") are preserved verbatim as text instead of being coerced
# to a floating point number by the COM Value setter.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.754.56'
$ws.Range("E2").Value = '  -1.17%  '
$ws.Range("D3").Value = '3.495.10'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '602.15'
$ws.Range("E5").Value = '  -1.13%  '
$ws.Range("D6").Value = '147.52'
$ws.Range("E6").Value = '  -3.19%  '
$ws.Range("D7").Value = '3.493.65'
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -1.74%  '
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("D11").Value = '7.80'
$ws.Range("E11").Value = '  +2.25%  '
$ws.Range("E12").Value = '  -2.31%  '
$ws.Range("E13").Value = '  -1.80%  '
$ws.Range("D14").Value = '4.086.48'
$ws.Range("D15").Value = '31.23'
$ws.Range("E15").Value = '  -4.56%  '
$ws.Range("D16").Value = '3.487.64'
$ws.Range("E16").Value = '  -0.63%  '
$ws.Range("D17").Value = '66.765.04'
$ws.Range("E17").Value = '  -1.15%  '
$ws.Range("E18").Value = '  -0.77%  '
$ws.Range("D19").Value = '10.50'
$ws.Range("E19").Value = '  +6.39%  '
$ws.Range("E20").Value = '  -3.00%  '
$ws.Range("E21").Value = '  -1.74%  '
$ws.Range("D22").Value = '433.45'
$ws.Range("E22").Value = '  -3.39%  '
$ws.Range("E23").Value = '  -4.32%  '
$ws.Range("D24").Value = '79.75'
$ws.Range("E24").Value = '  +1.96%  '
$ws.Range("D25").Value = '3.632.30'
$ws.Range("E25").Value = '  -0.52%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E28").Value = '  -7.01%  '
$ws.Range("D29").Value = '9.78'
$ws.Range("E29").Value = '  -3.36%  '
$ws.Range("D30").Value = '8.21'
$ws.Range("E30").Value = '  -7.28%  '
$ws.Range("D31").Value = '2.50'
$ws.Range("E31").Value = '  -1.01%  '
$ws.Range("E32").Value = '  -4.39%  '
$ws.Range("D33").Value = '0.998'
$ws.Range("E33").Value = '  -0.37%  '
$ws.Range("E34").Value = '  -1.92%  '
$ws.Range("D35").Value = '25.28'
$ws.Range("E35").Value = '  -1.96%  '
$ws.Range("D36").Value = '3.488.86'
$ws.Range("E36").Value = '  -0.42%  '
$ws.Range("E37").Value = '  -4.94%  '
$ws.Range("D38").Value = '5.87'
$ws.Range("E38").Value = '  -5.41%  '
$ws.Range("D39").Value = '7.98'
$ws.Range("E39").Value = '  -1.03%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("E42").Value = '  -1.20%  '
$ws.Range("D43").Value = '169.66'
$ws.Range("E43").Value = '  -2.58%  '
$ws.Range("E44").Value = '  -9.22%  '
$ws.Range("E45").Value = '  -1.55%  '
$ws.Range("E46").Value = '  +1.46%  '
$ws.Range("D47").Value = '28.84'
$ws.Range("E47").Value = '  -5.20%  '
$ws.Range("E48").Value = '  -2.12%  '
$ws.Range("D49").Value = '1.31'
$ws.Range("E49").Value = '  +0.50%  '
$ws.Range("E50").Value = '  -3.12%  '
$ws.Range("E51").Value = '  -4.89%  '

# Remove the temporary text format so the cells end up with the same
# (default/general) style they started with.
$ws.Range("D2:D51").ClearFormats()
